$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''302.49'
$ws.Range("E2").Value = '''-3.56%'
$ws.Range("G2").Value = '''12'
$ws.Range("D3").Value = '''35.25'
$ws.Range("E3").Value = '''-0.10%'
$ws.Range("G3").Value = '''12'
$ws.Range("D4").Value = '''5.038'
$ws.Range("E4").Value = '''-1.75%'
$ws.Range("G4").Value = '''12'
$ws.Range("D5").Value = '''0.07988'
$ws.Range("E5").Value = '''-1.71%'
$ws.Range("G5").Value = '''12'
$ws.Range("D6").Value = '''1.939'
$ws.Range("E6").Value = '''-8.77%'
$ws.Range("G6").Value = '''12'
$ws.Range("D7").Value = '''7.802'
$ws.Range("E7").Value = '''-2.00%'
$ws.Range("G7").Value = '''12'
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").Value = '''2.921'
$ws.Range("E8").Value = '''0.07%'
$ws.Range("G8").Value = '''12'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = '''0.9216'
$ws.Range("E9").Value = '''-0.69%'
$ws.Range("G9").Value = '''12'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '''0.1236'
$ws.Range("E10").Value = '''21.25%'
$ws.Range("G10").Value = '''12'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '''0.1853'
$ws.Range("E11").Value = '''-0.55%'
$ws.Range("G11").Value = '''12'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.09667'
$ws.Range("E12").Value = '''6.49%'
$ws.Range("G12").Value = '''12'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.03524'
$ws.Range("E13").Value = '''-1.93%'
$ws.Range("G13").Value = '''12'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.09858'
$ws.Range("E14").Value = '''-0.59%'
$ws.Range("G14").Value = '''12'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001398'
$ws.Range("E15").Value = '''-2.37%'
$ws.Range("G15").Value = '''12'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '''0.005818'
$ws.Range("E16").Value = '''2.04%'
$ws.Range("G16").Value = '''12'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '''3.502'
$ws.Range("E17").Value = '''0.90%'
$ws.Range("G17").Value = '''12'
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").Value = '''4.046'
$ws.Range("E18").Value = '''-2.50%'
$ws.Range("G18").Value = '''12'
$ws.Range("D19").Value = '''0.3401'
$ws.Range("E19").Value = '''-0.30%'
$ws.Range("G19").Value = '''12'
$ws.Range("D20").Value = '''0.1290'
$ws.Range("E20").Value = '''-3.83%'
$ws.Range("G20").Value = '''12'
$ws.Range("D21").Value = '''5.023'
$ws.Range("E21").Value = '''-1.40%'
$ws.Range("G21").Value = '''12'
$ws.Range("D22").Value = '''0.2394'
$ws.Range("E22").Value = '''8.11%'
$ws.Range("G22").Value = '''12'
$ws.Range("D23").Value = '''0.04484'
$ws.Range("E23").Value = '''-1.03%'
$ws.Range("G23").Value = '''12'
$ws.Range("E24").Value = '''-2.75%'
$ws.Range("G24").Value = '''12'
$ws.Range("D25").Value = '''0.004783'
$ws.Range("E25").Value = '''1.86%'
$ws.Range("G25").Value = '''12'
$ws.Range("D26").Value = '''0.0001248'
$ws.Range("E26").Value = '''-0.28%'
$ws.Range("G26").Value = '''12'
$ws.Range("D27").Value = '''0.0002996'
$ws.Range("E27").Value = '''-33.54%'
$ws.Range("G27").Value = '''12'
$ws.Range("G28").Value = '''12'
$ws.Range("G29").Value = '''12'
$ws.Range("G30").Value = '''12'
$ws.Range("G31").Value = '''12'
$ws.Range("G32").Value = '''12'
$ws.Range("G33").Value = '''12'
$ws.Range("G34").Value = '''12'
$ws.Range("G35").Value = '''12'
$ws.Range("G36").Value = '''12'
$ws.Range("G37").Value = '''12'
$ws.Range("G38").Value = '''12'
$ws.Range("D39").Value = '''0.01896'
$ws.Range("E39").Value = '''-3.42%'
$ws.Range("G39").Value = '''12'
$ws.Range("D40").Value = '''0.04729'
$ws.Range("E40").Value = '''-2.78%'
$ws.Range("G40").Value = '''12'
$ws.Range("D41").Value = '''0.007462'
$ws.Range("E41").Value = '''-2.52%'
$ws.Range("G41").Value = '''12'
$ws.Range("D42").Value = '''0.009885'
$ws.Range("E42").Value = '''25.98%'
$ws.Range("G42").Value = '''12'
$ws.Range("D43").Value = '''0.1327'
$ws.Range("E43").Value = '''-4.50%'
$ws.Range("G43").Value = '''12'
$ws.Range("D44").Value = '''0.002106'
$ws.Range("E44").Value = '''-0.30%'
$ws.Range("G44").Value = '''12'
$ws.Range("D45").Value = '''0.01084'
$ws.Range("E45").Value = '''-7.73%'
$ws.Range("G45").Value = '''12'
$ws.Range("D46").Value = '''0.00006240'
$ws.Range("E46").Value = '''-6.89%'
$ws.Range("G46").Value = '''12'
$ws.Range("D47").Value = '''0.00000000749'
$ws.Range("E47").Value = '''-0.30%'
$ws.Range("G47").Value = '''12'
$ws.Range("E48").Value = '''68.24%'
$ws.Range("G48").Value = '''12'
$ws.Range("E49").Value = '''-12.53%'
$ws.Range("G49").Value = '''12'
$ws.Range("D50").Value = '''0.00002097'
$ws.Range("E50").Value = '''-0.30%'
$ws.Range("G50").Value = '''12'
$ws.Range("D51").Value = '''0.0001997'
$ws.Range("E51").Value = '''-0.30%'
$ws.Range("G51").Value = '''12'
